$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# ---------------------------------------------------------------------------
# 1) "Weekly Quantity": append 3 new weekly rows
# ---------------------------------------------------------------------------
$wsWeekly.Range("A2").Copy()
$wsWeekly.Range("A3:A5").PasteSpecial(-4122)

$weeklyRows = @(
  @(3, 45662.99999999999, 3),
  @(4, 45669.99999999999, 6),
  @(5, 45683.99999999999, 1)
)
foreach ($row in $weeklyRows) {
  $r = $row[0]
  $wsWeekly.Cells.Item($r, 1).Value = $row[1]
  $wsWeekly.Cells.Item($r, 2).Value = $row[2]
}

# ---------------------------------------------------------------------------
# 2) "Monthly Trend": append 1 new monthly row
# ---------------------------------------------------------------------------
$wsMonthly.Range("A2").Copy()
$wsMonthly.Range("A3:A3").PasteSpecial(-4122)

$wsMonthly.Cells.Item(3, 1).Value = 45688.99999999999
$wsMonthly.Cells.Item(3, 2).Value = 10

# ---------------------------------------------------------------------------
# 3) New sheet "PO Forecast" at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the page margins used on the other sheets (0.75in/1in/0.5in)
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Header formatting (bold / bordered / centered), matching the other sheets
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:B1").PasteSpecial(-4122)
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"

# Date-column formatting for rows 2-13
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A13").PasteSpecial(-4122)

$forecastRows = @(
  @(2,  45613.99999999999, 115),
  @(3,  45662.99999999999, 23),
  @(4,  45669.99999999999, 10),
  @(5,  45683.99999999999, 0),
  @(6,  45690.99999999999, 0),
  @(7,  45697.99999999999, 0),
  @(8,  45704.99999999999, 0),
  @(9,  45711.99999999999, 0),
  @(10, 45718.99999999999, 0),
  @(11, 45725.99999999999, 0),
  @(12, 45732.99999999999, 0),
  @(13, 45739.99999999999, 0)
)
foreach ($row in $forecastRows) {
  $r = $row[0]
  $wsForecast.Cells.Item($r, 1).Value = $row[1]
  $wsForecast.Cells.Item($r, 2).Value = $row[2]
}

Write-Host "PO Forecast model applied"
